$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cell, [string]$value, [bool]$forceText = $false)
    if ($forceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}

# Row 2
Set-CellText $ws.Cells.Item(2, 4) "60.278.77" $false
Set-CellText $ws.Cells.Item(2, 5) "  +1.25%  " $false

# Row 3
Set-CellText $ws.Cells.Item(3, 4) "2.585.57" $false
Set-CellText $ws.Cells.Item(3, 5) "  +0.76%  " $false

# Row 4
Set-CellText $ws.Cells.Item(4, 5) "  -0.17%  " $false

# Row 5
Set-CellText $ws.Cells.Item(5, 4) "505.86" $true
Set-CellText $ws.Cells.Item(5, 5) "  -0.01%  " $false

# Row 6
Set-CellText $ws.Cells.Item(6, 4) "152.50" $true
Set-CellText $ws.Cells.Item(6, 5) "  -4.30%  " $false

# Row 7
Set-CellText $ws.Cells.Item(7, 5) "  +0.41%  " $false

# Row 8
Set-CellText $ws.Cells.Item(8, 5) "  -6.30%  " $false

# Row 9
Set-CellText $ws.Cells.Item(9, 4) "2.589.85" $false
Set-CellText $ws.Cells.Item(9, 5) "  -0.01%  " $false

# Row 10
Set-CellText $ws.Cells.Item(10, 4) "6.70" $true
Set-CellText $ws.Cells.Item(10, 5) "  +7.65%  " $false

# Row 11
Set-CellText $ws.Cells.Item(11, 4) "0.103" $true
Set-CellText $ws.Cells.Item(11, 5) "  -1.49%  " $false

# Row 12
Set-CellText $ws.Cells.Item(12, 5) "  +0.33%  " $false

# Row 13
Set-CellText $ws.Cells.Item(13, 5) "  +0.67%  " $false

# Row 14
Set-CellText $ws.Cells.Item(14, 4) "3.037.71" $false
Set-CellText $ws.Cells.Item(14, 5) "  +2.01%  " $false

# Row 15
Set-CellText $ws.Cells.Item(15, 4) "60.232.94" $false
Set-CellText $ws.Cells.Item(15, 5) "  +1.72%  " $false

# Row 16
Set-CellText $ws.Cells.Item(16, 4) "21.54" $true
Set-CellText $ws.Cells.Item(16, 5) "  -2.18%  " $false

# Row 17
Set-CellText $ws.Cells.Item(17, 5) "  +0.78%  " $false

# Row 18
Set-CellText $ws.Cells.Item(18, 4) "2.587.78" $false
Set-CellText $ws.Cells.Item(18, 5) "  +0.61%  " $false

# Row 19
Set-CellText $ws.Cells.Item(19, 4) "4.80" $true
Set-CellText $ws.Cells.Item(19, 5) "  +1.00%  " $false

# Row 20
Set-CellText $ws.Cells.Item(20, 4) "346.27" $true
Set-CellText $ws.Cells.Item(20, 5) "  +3.63%  " $false

# Row 21
Set-CellText $ws.Cells.Item(21, 5) "  -0.43%  " $false

# Row 22
Set-CellText $ws.Cells.Item(22, 4) "6.08" $true
Set-CellText $ws.Cells.Item(22, 5) "  -0.05%  " $false

# Row 23
Set-CellText $ws.Cells.Item(23, 5) "  -0.61%  " $false

# Row 24
Set-CellText $ws.Cells.Item(24, 4) "60.25" $true
Set-CellText $ws.Cells.Item(24, 5) "  +0.22%  " $false

# Row 25
Set-CellText $ws.Cells.Item(25, 4) "0.418" $true
Set-CellText $ws.Cells.Item(25, 5) "  -0.01%  " $false

# Row 26
Set-CellText $ws.Cells.Item(26, 2) "Kaspa" $false
Set-CellText $ws.Cells.Item(26, 3) "https://coinranking.com/coin/V8GxkwWow+kaspa-kas" $false
Set-CellText $ws.Cells.Item(26, 4) "0.165" $true
Set-CellText $ws.Cells.Item(26, 5) "  -1.14%  " $false

# Row 27
Set-CellText $ws.Cells.Item(27, 2) "WrappedeETH" $false
Set-CellText $ws.Cells.Item(27, 3) "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth" $false
Set-CellText $ws.Cells.Item(27, 4) "2.695.98" $false
Set-CellText $ws.Cells.Item(27, 5) "  +1.98%  " $false

# Row 28
Set-CellText $ws.Cells.Item(28, 4) "1.00" $true

# Row 29
Set-CellText $ws.Cells.Item(29, 5) "  +1.29%  " $false

# Row 30
Set-CellText $ws.Cells.Item(30, 5) "  -1.64%  " $false

# Row 31
Set-CellText $ws.Cells.Item(31, 5) "  +0.27%  " $false

# Row 32
Set-CellText $ws.Cells.Item(32, 4) "155.00" $true
Set-CellText $ws.Cells.Item(32, 5) "  -0.18%  " $false

# Row 33
Set-CellText $ws.Cells.Item(33, 4) "19.23" $true
Set-CellText $ws.Cells.Item(33, 5) "  -1.43%  " $false

# Row 34
Set-CellText $ws.Cells.Item(34, 5) "  -0.94%  " $false

# Row 35
Set-CellText $ws.Cells.Item(35, 5) "  +3.62%  " $false

# Row 36
Set-CellText $ws.Cells.Item(36, 4) "3.99" $true
Set-CellText $ws.Cells.Item(36, 5) "  +1.67%  " $false

# Row 37
Set-CellText $ws.Cells.Item(37, 5) "  -0.75%  " $false

# Row 38
Set-CellText $ws.Cells.Item(38, 4) "0.860" $true
Set-CellText $ws.Cells.Item(38, 5) "  +20.24%  " $false

# Row 39
Set-CellText $ws.Cells.Item(39, 2) "Filecoin" $false
Set-CellText $ws.Cells.Item(39, 3) "https://coinranking.com/coin/ymQub4fuB+filecoin-fil" $false
Set-CellText $ws.Cells.Item(39, 4) "3.77" $true
Set-CellText $ws.Cells.Item(39, 5) "  +0.28%  " $false

# Row 40
Set-CellText $ws.Cells.Item(40, 2) "Fetch.AI" $false
Set-CellText $ws.Cells.Item(40, 3) "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet" $false
Set-CellText $ws.Cells.Item(40, 4) "0.843" $true
Set-CellText $ws.Cells.Item(40, 5) "  -2.25%  " $false

# Row 41
Set-CellText $ws.Cells.Item(41, 4) "299.23" $true
Set-CellText $ws.Cells.Item(41, 5) "  +2.83%  " $false

# Row 42
Set-CellText $ws.Cells.Item(42, 4) "35.72" $true
Set-CellText $ws.Cells.Item(42, 5) "  +2.62%  " $false

# Row 43
Set-CellText $ws.Cells.Item(43, 5) "  +0.54%  " $false

# Row 44
Set-CellText $ws.Cells.Item(44, 4) "0.617" $true
Set-CellText $ws.Cells.Item(44, 5) "  -1.65%  " $false

# Row 45
Set-CellText $ws.Cells.Item(45, 5) "  -2.09%  " $false

# Row 46
Set-CellText $ws.Cells.Item(46, 4) "0.0558" $true
Set-CellText $ws.Cells.Item(46, 5) "  -0.25%  " $false

# Row 47
Set-CellText $ws.Cells.Item(47, 5) "  +0.02%  " $false

# Row 48
Set-CellText $ws.Cells.Item(48, 4) "19.70" $true
Set-CellText $ws.Cells.Item(48, 5) "  +2.38%  " $false

# Row 49
Set-CellText $ws.Cells.Item(49, 4) "4.85" $true
Set-CellText $ws.Cells.Item(49, 5) "  +1.44%  " $false

# Row 50
Set-CellText $ws.Cells.Item(50, 5) "  -2.53%  " $false

# Row 51
Set-CellText $ws.Cells.Item(51, 4) "2.014.11" $false

